# Update the generated-experiment numeric/expression values to the new
# "alpha_zero" scenario figures (todos no convexos menos el 5to).
#
# Every one of these cells is stored in the workbook as a *text* shared
# string (t="s"), even the ones that look like plain numbers (e.g.
# "-8.05"). Assigning a numeric-looking string straight to .Value makes
# Excel coerce it into a real number cell, so we briefly force the
# target ranges to Text number format before writing the values, then
# clear the formatting again so the cells end up back at the sheet's
# default (General) style - matching the original file, just with new
# text contents.

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

$wsFollower.Range("A2:F5").NumberFormat = "@"

$wsFollower.Range("A2").Value = "8.05 - y"
$wsFollower.Range("B2").Value = "-8.05"
$wsFollower.Range("D2").Value = "0.13"
$wsFollower.Range("E2").Value = "6.0"
$wsFollower.Range("F2").Value = "3.7"

$wsFollower.Range("A3").Value = "-1.950000000000001 - x + y"
$wsFollower.Range("B3").Value = "-1.049999999999999"
$wsFollower.Range("D3").Value = "0.6"
$wsFollower.Range("E3").Value = "7.9"
$wsFollower.Range("F3").Value = "2.1"

$wsFollower.Range("A4").Value = "-22.200000000000003 + x + 2y"
$wsFollower.Range("B4").Value = "10.200000000000001"
$wsFollower.Range("D4").Value = "0.62"
$wsFollower.Range("E4").Value = "9.8"
$wsFollower.Range("F4").Value = "0.2"

$wsFollower.Range("A5").Value = "-16.949999999999996 + 4x - y"
$wsFollower.Range("B5").Value = "4.349999999999998"
$wsFollower.Range("D5").Value = "0.32"
$wsFollower.Range("E5").Value = "1.7999999999999998"
$wsFollower.Range("F5").Value = "9.0"

$wsFollower.Range("A2:F5").ClearFormats()

# --- Punto_modificado ---------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")

$wsPunto.Range("A2:B2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "6.1"
$wsPunto.Range("B2").Value = "8.05"
$wsPunto.Range("A2:B2").ClearFormats()

# --- Vector_bf ------------------------------------------------------------
# NOTE: worksheet-name lookups (and PowerShell variable names) are
# case-insensitive, and this workbook has both a "Vector_bf" sheet and a
# "Vector_BF" sheet whose names differ only by case. Looking either one up
# by name would resolve to whichever of the two comes first, so both are
# addressed here by their (1-based) tab position instead.
$wsVecLower = $wb.Worksheets.Item(5)   # "Vector_bf"

$wsVecLower.Range("A2").NumberFormat = "@"
$wsVecLower.Range("A2").Value = "-2.3899999999999997"
$wsVecLower.Range("A2").ClearFormats()

# --- Vector_BF --------------------------------------------------------------
$wsVecUpper = $wb.Worksheets.Item(6)   # "Vector_BF"

$wsVecUpper.Range("A2:A3").NumberFormat = "@"
$wsVecUpper.Range("A2").Value = "-8.1"
$wsVecUpper.Range("A3").Value = "-16.7"
$wsVecUpper.Range("A2:A3").ClearFormats()
